$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new data rows (58, 59) ------------------------------------
# Copy formatting from the last existing row (57) down into the new rows so
# that columns A/B pick up the same "data row" style already used throughout
# the table, then overwrite the values.
$ws.Range("A57:C57").Copy() | Out-Null
$ws.Range("A58:C59").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A58").Value = "GENO"
$ws.Range("B58").Value = "GENIO"
$ws.Range("C58").Value = "Commercial"

$ws.Range("A59").Value = "LOGN"
$ws.Range("B59").Value = "LOGAN"
$ws.Range("C59").Value = "Personal"

# --- Apply a thin border down the whole Segment column (C2:C59) ------------
$ws.Range("C2:C59").Borders.LineStyle = 1

# --- Move the active selection (was C6, now B6) -----------------------------
$ws.Range("B6").Select() | Out-Null

# --- Register the (hidden) filter-database name for the table range --------
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$59")
$nm.Visible = $false
